$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "ba"
$ws.Range("J2").Value = "Appreciation"
$ws.Range("I33").Value = "sd"
$ws.Range("J33").Value = "Statement-non-opinion"
$ws.Range("I41").Value = "sv"
$ws.Range("J41").Value = "Statement-opinion"
$ws.Range("I57").Value = "sd"
$ws.Range("J57").Value = "Statement-non-opinion"
$ws.Range("I60").Value = "sv"
$ws.Range("J60").Value = "Statement-opinion"
$ws.Range("I78").Value = "aa"
$ws.Range("J78").Value = "Agree/Accept"
$ws.Range("I121").Value = "b"
$ws.Range("J121").Value = "Acknowledge (Backchannel)"
$ws.Range("I125").Value = "aa"
$ws.Range("J125").Value = "Agree/Accept"
$ws.Range("I132").Value = "sd"
$ws.Range("J132").Value = "Statement-non-opinion"
$ws.Range("I147").Value = "b"
$ws.Range("J147").Value = "Acknowledge (Backchannel)"
$ws.Range("I149").Value = "sv"
$ws.Range("J149").Value = "Statement-opinion"
$ws.Range("I162").Value = "ba"
$ws.Range("J162").Value = "Appreciation"
$ws.Range("I171").Value = "%"
$ws.Range("J171").Value = "Uninterpretable"
$ws.Range("I173").Value = "aa"
$ws.Range("J173").Value = "Agree/Accept"
$ws.Range("I182").Value = "%"
$ws.Range("J182").Value = "Uninterpretable"
$ws.Range("I184").Value = "sv"
$ws.Range("J184").Value = "Statement-opinion"
$ws.Range("I185").Value = "aa"
$ws.Range("J185").Value = "Agree/Accept"
$ws.Range("I188").Value = "aa"
$ws.Range("J188").Value = "Agree/Accept"
$ws.Range("I221").Value = "b"
$ws.Range("J221").Value = "Acknowledge (Backchannel)"
$ws.Range("I223").Value = "sd"
$ws.Range("J223").Value = "Statement-non-opinion"
$ws.Range("I231").Value = "aa"
$ws.Range("J231").Value = "Agree/Accept"
$ws.Range("I237").Value = "b"
$ws.Range("J237").Value = "Acknowledge (Backchannel)"
$ws.Range("I252").Value = "sd"
$ws.Range("J252").Value = "Statement-non-opinion"
$ws.Range("I257").Value = "sd"
$ws.Range("J257").Value = "Statement-non-opinion"
$ws.Range("I262").Value = "%"
$ws.Range("J262").Value = "Uninterpretable"
$ws.Range("I274").Value = "sd"
$ws.Range("J274").Value = "Statement-non-opinion"
$ws.Range("I283").Value = "%"
$ws.Range("J283").Value = "Uninterpretable"
$ws.Range("I284").Value = "sv"
$ws.Range("J284").Value = "Statement-opinion"
$ws.Range("I297").Value = "aa"
$ws.Range("J297").Value = "Agree/Accept"
$ws.Range("I301").Value = "ba"
$ws.Range("J301").Value = "Appreciation"
$ws.Range("I302").Value = "sd"
$ws.Range("J302").Value = "Statement-non-opinion"
$ws.Range("I303").Value = "sd"
$ws.Range("J303").Value = "Statement-non-opinion"
$ws.Range("I348").Value = "sv"
$ws.Range("J348").Value = "Statement-opinion"
$ws.Range("I373").Value = "%"
$ws.Range("J373").Value = "Uninterpretable"
$ws.Range("I387").Value = "sd"
$ws.Range("J387").Value = "Statement-non-opinion"
$ws.Range("I401").Value = "aa"
$ws.Range("J401").Value = "Agree/Accept"
$ws.Range("I405").Value = "%"
$ws.Range("J405").Value = "Uninterpretable"
$ws.Range("I407").Value = "sd"
$ws.Range("J407").Value = "Statement-non-opinion"
$ws.Range("I417").Value = "sd"
$ws.Range("J417").Value = "Statement-non-opinion"
$ws.Range("I423").Value = "sd"
$ws.Range("J423").Value = "Statement-non-opinion"
$ws.Range("I429").Value = "sd"
$ws.Range("J429").Value = "Statement-non-opinion"
$ws.Range("I436").Value = "sv"
$ws.Range("J436").Value = "Statement-opinion"
$ws.Range("I441").Value = "%"
$ws.Range("J441").Value = "Uninterpretable"
$ws.Range("I443").Value = "sv"
$ws.Range("J443").Value = "Statement-opinion"
$ws.Range("I448").Value = "sv"
$ws.Range("J448").Value = "Statement-opinion"
$ws.Range("I455").Value = "ba"
$ws.Range("J455").Value = "Appreciation"
$ws.Range("I461").Value = "%"
$ws.Range("J461").Value = "Uninterpretable"
$ws.Range("I463").Value = "aa"
$ws.Range("J463").Value = "Agree/Accept"
$ws.Range("I482").Value = "b"
$ws.Range("J482").Value = "Acknowledge (Backchannel)"
$ws.Range("I484").Value = "aa"
$ws.Range("J484").Value = "Agree/Accept"
$ws.Range("I485").Value = "b"
$ws.Range("J485").Value = "Acknowledge (Backchannel)"
$ws.Range("I487").Value = "sd"
$ws.Range("J487").Value = "Statement-non-opinion"
$ws.Range("I498").Value = "sv"
$ws.Range("J498").Value = "Statement-opinion"
$ws.Range("I499").Value = "sd"
$ws.Range("J499").Value = "Statement-non-opinion"
$ws.Range("I511").Value = "sv"
$ws.Range("J511").Value = "Statement-opinion"
$ws.Range("I516").Value = "sd"
$ws.Range("J516").Value = "Statement-non-opinion"
$ws.Range("I522").Value = "sd"
$ws.Range("J522").Value = "Statement-non-opinion"
$ws.Range("I526").Value = "sv"
$ws.Range("J526").Value = "Statement-opinion"
$ws.Range("I528").Value = "ba"
$ws.Range("J528").Value = "Appreciation"
$ws.Range("I532").Value = "b"
$ws.Range("J532").Value = "Acknowledge (Backchannel)"
$ws.Range("I543").Value = "ba"
$ws.Range("J543").Value = "Appreciation"
$ws.Range("I554").Value = "aa"
$ws.Range("J554").Value = "Agree/Accept"
$ws.Range("I555").Value = "sd"
$ws.Range("J555").Value = "Statement-non-opinion"
$ws.Range("I570").Value = "aa"
$ws.Range("J570").Value = "Agree/Accept"
$ws.Range("I573").Value = "sd"
$ws.Range("J573").Value = "Statement-non-opinion"
$ws.Range("I589").Value = "b"
$ws.Range("J589").Value = "Acknowledge (Backchannel)"
$ws.Range("I591").Value = "sd"
$ws.Range("J591").Value = "Statement-non-opinion"
